$d = $word.ActiveDocument

$replacements = @(
    @{old = "789÷7=112, 5"; new = "864÷6=144, 0"},
    @{old = "909÷2=454, 1"; new = "252÷9=28, 0"},
    @{old = "191÷9=21, 2"; new = "826÷9=91, 7"},
    @{old = "249÷3=83, 0"; new = "390÷5=78, 0"},
    @{old = "753÷7=107, 4"; new = "741÷7=105, 6"},
    @{old = "215÷8=26, 7"; new = "532÷3=177, 1"},
    @{old = "331÷6=55, 1"; new = "803÷8=100, 3"},
    @{old = "403÷5=80, 3"; new = "945÷7=135, 0"},
    @{old = "247÷6=41, 1"; new = "356÷2=178, 0"},
    @{old = "677÷4=169, 1"; new = "930÷5=186, 0"},
    @{old = "258÷8=32, 2"; new = "127÷9=14, 1"},
    @{old = "605÷3=201, 2"; new = "746÷7=106, 4"},
    @{old = "937÷2=468, 1"; new = "848÷2=424, 0"},
    @{old = "885÷5=177, 0"; new = "736÷6=122, 4"},
    @{old = "542÷3=180, 2"; new = "712÷4=178, 0"},
    @{old = "610÷3=203, 1"; new = "655÷6=109, 1"},
    @{old = "903÷2=451, 1"; new = "594÷3=198, 0"},
    @{old = "355÷9=39, 4"; new = "351÷7=50, 1"},
    @{old = "627÷2=313, 1"; new = "519÷7=74, 1"},
    @{old = "709÷8=88, 5"; new = "782÷3=260, 2"},
    @{old = "942÷6=157, 0"; new = "114÷2=57, 0"},
    @{old = "489÷4=122, 1"; new = "537÷7=76, 5"},
    @{old = "778÷9=86, 4"; new = "819÷4=204, 3"},
    @{old = "254÷8=31, 6"; new = "231÷3=77, 0"},
    @{old = "234÷4=58, 2"; new = "873÷6=145, 3"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
